$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert two new columns (mirrors the structural change in the diff) ---
# Insert a new column before column O (15)
$ws.Columns.Item(15).Insert()
# Insert a new column before the (new) column R (18) -- this was old column Q's new position
$ws.Columns.Item(18).Insert()

# --- Step 2: fix up row 1 (figure captions) ---
# Old O1 caption shifted to P1 on the first insert; clear it, the new O1 gets a *new* caption.
$ws.Range("P1").ClearContents()
$ws.Range("O1").Value = "Figure 1: Annual Growth Rates for Housing and Wages in Australia (1990–2024). Source: ABS."
$ws.Range("O1").Style = $ws.Range("Q1").Style
# S1 already holds the old "Fig 2" caption (shifted via both inserts); just refresh its text.
$ws.Range("S1").Value = "Fig 2,  Indexed Median House Price vs. Wage Growth since 1975. Source: ABS."

# --- Step 3: header row 2, new helper "Year" columns ---
$ws.Range("O2").Value = "Year"
$ws.Range("R2").Value = "Year"
$ws.Range("R2").Style = $ws.Range("O2").Style

# --- Step 4: fill literal Year values down columns O and R for every data row ---
for ($r = 3; $r -le 53; $r++) {
    $yearVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 15).Value = $yearVal   # column O = A (Year)
    $ws.Cells.Item($r, 18).Value = $yearVal   # column R = A (Year)
    $ws.Cells.Item($r, 18).Style = $ws.Cells.Item($r, 15).Style
}
